$wb = $excel.ActiveWorkbook

# --- "ops" sheet: fix the "sumbit-data" typo -> "submit-data" and repoint
# both operation-definition hyperlink URLs away from the (now invalid)
# davinci-deqm STU3 paths to the canonical hl7.org/fhir OperationDefinition
# URLs. Also drop the stray hyperlink that had been left on B2 (it still
# pointed at "collect-data").
$ops = $wb.Worksheets.Item("ops")

$ops.Range("B2").Hyperlinks.Delete() | Out-Null

$ops.Range("A2").Value = "submit-data"
$ops.Range("B2").Value = "http://hl7.org/fhir/OperationDefinition/Measure-submit-data"

$ops.Range("A3").Value = "data-requirements"
$ops.Range("B3").Value = "http://hl7.org/fhir/OperationDefinition/Measure-data-requirements"

# --- Active sheet/selection bookkeeping: the workbook was last saved with
# "ops" as the active tab (selection resting on B2), whereas it had
# previously been saved with "meta" active (selection resting on B8).
$meta = $wb.Worksheets.Item("meta")
$meta.Range("B8").Select() | Out-Null

$ops.Activate() | Out-Null
$ops.Range("B2").Select() | Out-Null
